# edit.ps1
# Applies the "feat: add 2022-Q3 data" change:
#  1. Updates the "总计" (totals) sheet: inserts a new first data row for 2022-Q3
#     and keeps the existing quarters, shifting them down by one row.
#  2. Inserts a brand new worksheet named "2022-Q3" (positioned right after
#     "总计", before the former first quarter sheet "2022-Q2") containing the
#     per-fund holding breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: update the "总计" summary sheet
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Extend the styled index column (A) down into the new row 6 by copying the
# formatting that is already used by the existing index cells.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q3", 12, 0.46),
    @("2022-Q2", 2, 0.15),
    @("2021-Q4", 15, 2.52),
    @("2021-Q3", 25, 4.58),
    @("2021-Q2", 1, 0.01)
)

$r = 2
foreach ($row in $totalRows) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Part 2: insert the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
# The former first-quarter sheet ("2022-Q2") already has the exact layout
# (headers, styles, column types) that the new quarter sheet needs, so
# duplicate it and place the copy immediately before it; this keeps every
# style/format identical to the sibling quarter sheets.
$existingQ2 = $wb.Worksheets.Item(2)
$existingQ2.Copy($existingQ2, $null)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template only had 2 data rows; extend the styled rows down to 13
# (12 data rows) by copying the formatting of row 2.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A4:H13").PasteSpecial(-4122)

$fundRows = @(
    @("010114","华宝新兴成长混合","3.18","80.50","3.73","0.1186",9),
    @("014600","博时回报严选混合A","0.92","92.53","9.43","0.0868",1),
    @("014232","博时专精特新主题混合A","3.14","81.93","2.48","0.0779",1),
    @("014233","博时专精特新主题混合C","2.69","81.93","2.48","0.0667",1),
    @("011927","博时汇誉回报混合A","1.04","80.45","5.39","0.0561",5),
    @("006813","博时汇悦回报混合","0.79","79.83","3.78","0.0299",5),
    @("011928","博时汇誉回报混合C","0.12","80.45","5.39","0.0065",5),
    @("501002","长信价值优选混合","0.39","93.83","1.61","0.0063",4),
    @("014601","博时回报严选混合C","0.04","92.53","9.43","0.0038",1),
    @("012415","德邦上证G60综指增强A","0.09","92.83","3.09","0.0028",3),
    @("000822","东海美丽中国灵活配置混合","0.26","25.44","0.89","0.0023",6),
    @("012416","德邦上证G60综指增强C","0.01","92.83","3.09","0.0003",3)

)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
